$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Strikethrough the "Hans XP-Archiv" row (row 9) - device being replaced.
$ws.Range("A9:K9").Font.Strikethrough = $true
$ws.Range("L9").Font.Strikethrough = $true

# 2. Insert a new row for "Lae12 Hans Win10" before the current row 11
#    (Excel naturally inherits the formatting of the row above on insert).
$ws.Rows("11:11").Insert()

$ws.Range("G11").Value = "Lae12 Hans Win10"
$ws.Range("A11").Value = "Hans Win10"
$ws.Range("B11").Value = 80
$ws.Range("C11").Value = 917
$ws.Range("D11").Formula = "=B11/C11"
$ws.Range("F11").Value = 140
$ws.Range("L11").Value = 42597
$ws.Range("L11").NumberFormat = "m/d/yyyy"

# 3. Restore cursor/selection to roughly where the author left it.
$ws.Range("L12").Select()
